$d = $word.ActiveDocument
$d.Content.Find.Execute("Revision de contrato", $true, $false, $false, $false, $false,
                         $true, 1, $false, "revisión de contrato", 2)
